$wb = $excel.ActiveWorkbook

# Fix typo in sheet name: "busbar" -> "bus"
$wsBus = $wb.Worksheets.Item("busbar")
$wsBus.Name = "bus"

# Narrow the generator sheet's frozen-pane selection from R3:R50 down to R3
$wsGen = $wb.Worksheets.Item("generator")
$wsGen.Activate()
$wsGen.Range("R3").Select()

# Make the renamed "bus" sheet the active tab/selection (was "generator" before)
$wsBus.Activate()
$wsBus.Range("A4").Select()
